# repull data, push all data, mean calculation
# Update the dSF column (F) values for the jake irvin save-data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = 0
    4  = 2
    7  = 2
    8  = -1
    9  = 2
    10 = 1
    11 = -3
    12 = 2
    13 = -2
    14 = 6
    15 = -2
    16 = -1
    17 = -1
    19 = 1
    20 = -3
    21 = -1
    22 = -6
    23 = -2
    26 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
